$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '52.132.66'
$ws.Range("E2").Value = '  +0.04%  '
$ws.Range("D3").Value = '2.843.29'
$ws.Range("E3").Value = '  +1.81%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '361.30'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +5.36%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '113.88'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.40%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.570'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +4.79%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.604'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.98%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.72'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.90%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0863'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.75%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.132'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.23%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.99'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.67%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.77'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.72%  '
$ws.Range("D15").Value = '3.287.83'
$ws.Range("E15").Value = '  +1.64%  '
$ws.Range("D16").Value = '2.836.08'
$ws.Range("E16").Value = '  +0.73%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.915'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.08%  '
$ws.Range("D18").Value = '51.909.40'
$ws.Range("E18").Value = '  -0.25%  '
$ws.Range("E19").Value = '  +7.01%  '
$ws.Range("E20").Value = '  -1.64%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.61'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.97%  '
$ws.Range("D22").Value = '0.0₃0994'
$ws.Range("E22").Value = '  +0.32%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.29'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.06%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '267.16'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.82%  '
$ws.Range("E25").Value = '  -0.03%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '27.24'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.33%  '
$ws.Range("E27").Value = '  +0.05%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.44'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.60%  '
$ws.Range("E29").Value = '  +1.57%  '
$ws.Range("E30").Value = '  +5.20%  '
$ws.Range("E31").Value = '  -1.70%  '
$ws.Range("B32").Value = 'VeChain'
$ws.Range("C32").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0458'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +23.53%  '
$ws.Range("B33").Value = 'InjectiveProtocol'
$ws.Range("C33").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '34.17'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.55%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.90'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.92%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.38'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +8.09%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0844'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.44%  '
$ws.Range("E37").Value = '  -0.06%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.30'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.71%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.08'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.42%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.34'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.32%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '23.99'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.79%  '
$ws.Range("E42").Value = '  +2.18%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '128.08'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.14%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.56'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -7.02%  '
$ws.Range("E45").Value = '  -3.21%  '
$ws.Range("D46").Value = '2.123.44'
$ws.Range("E46").Value = '  +0.73%  '
$ws.Range("E47").Value = '  +1.61%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.01'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +9.71%  '
$ws.Range("E50").Value = '  +5.21%  '
$ws.Range("E51").Value = '  +1.20%  '
